$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number need to be forced to stay
# text (matching the original inlineStr/text cell type) by temporarily applying
# a Text number format, then restoring the default "Normal" style so no stray
# cell-level style attribute is left behind.
$textCells = @(
    "D5",
    "D6",
    "D7",
    "D9",
    "D10",
    "D13",
    "D14",
    "D15",
    "D19",
    "D20",
    "D22",
    "D26",
    "D27",
    "D29",
    "D30",
    "D32",
    "D33",
    "D34",
    "D37",
    "D38",
    "D40",
    "D41",
    "D43",
    "D44",
    "D47",
    "D50",
)
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = '38.126.54'
$ws.Range("E2").Value = '  +3.04%  '
$ws.Range("D3").Value = '2.061.12'
$ws.Range("E3").Value = '  +2.37%  '
$ws.Range("E4").Value = '  +0.35%  '
$ws.Range("D5").Value = '230.63'
$ws.Range("E5").Value = '  +2.37%  '
$ws.Range("D6").Value = '0.615'
$ws.Range("E6").Value = '  +1.59%  '
$ws.Range("D7").Value = '58.37'
$ws.Range("E7").Value = '  +6.53%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = '0.387'
$ws.Range("E9").Value = '  +2.93%  '
$ws.Range("D10").Value = '0.0809'
$ws.Range("E10").Value = '  +3.62%  '
$ws.Range("E11").Value = '  +0.70%  '
$ws.Range("D12").Value = '2.369.11'
$ws.Range("E12").Value = '  +2.72%  '
$ws.Range("D13").Value = '14.63'
$ws.Range("E13").Value = '  +3.78%  '
$ws.Range("D14").Value = '20.68'
$ws.Range("E14").Value = '  +2.60%  '
$ws.Range("D15").Value = '0.753'
$ws.Range("E15").Value = '  +1.92%  '
$ws.Range("E16").Value = '  +3.48%  '
$ws.Range("D17").Value = '2.066.76'
$ws.Range("E17").Value = '  +2.66%  '
$ws.Range("D18").Value = '38.019.20'
$ws.Range("E18").Value = '  +2.78%  '
$ws.Range("D19").Value = '6.14'
$ws.Range("E19").Value = '  -0.64%  '
$ws.Range("D20").Value = '69.95'
$ws.Range("E20").Value = '  +1.91%  '
$ws.Range("D21").Value = '0.0₃0831'
$ws.Range("E21").Value = '  +2.28%  '
$ws.Range("D22").Value = '224.92'
$ws.Range("E22").Value = '  +1.00%  '
$ws.Range("E24").Value = '  +1.15%  '
$ws.Range("E25").Value = '  +3.28%  '
$ws.Range("D26").Value = '9.30'
$ws.Range("E26").Value = '  +1.93%  '
$ws.Range("D27").Value = '166.28'
$ws.Range("E27").Value = '  +0.25%  '
$ws.Range("E28").Value = '  +7.16%  '
$ws.Range("D29").Value = '19.05'
$ws.Range("E29").Value = '  +1.95%  '
$ws.Range("D30").Value = '1.35'
$ws.Range("E30").Value = '  +0.18%  '
$ws.Range("E31").Value = '  +1.84%  '
$ws.Range("D32").Value = '4.55'
$ws.Range("E32").Value = '  +0.88%  '
$ws.Range("D33").Value = '4.62'
$ws.Range("E33").Value = '  +5.16%  '
$ws.Range("D34").Value = '0.0613'
$ws.Range("E34").Value = '  +0.36%  '
$ws.Range("E35").Value = '  +7.57%  '
$ws.Range("E36").Value = '  +0.61%  '
$ws.Range("D37").Value = '6.03'
$ws.Range("E37").Value = '  +14.54%  '
$ws.Range("D38").Value = '3.33'
$ws.Range("E38").Value = '  +5.52%  '
$ws.Range("E39").Value = '  -0.05%  '
$ws.Range("D40").Value = '0.0219'
$ws.Range("E40").Value = '  +2.06%  '
$ws.Range("D41").Value = '98.45'
$ws.Range("E41").Value = '  +3.86%  '
$ws.Range("D42").Value = '1.480.86'
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("D43").Value = '0.0943'
$ws.Range("E43").Value = '  +3.09%  '
$ws.Range("D44").Value = '16.78'
$ws.Range("E44").Value = '  +3.20%  '
$ws.Range("E45").Value = '  +3.57%  '
$ws.Range("E46").Value = '  +0.12%  '
$ws.Range("D47").Value = '4.10'
$ws.Range("E48").Value = '  +1.30%  '
$ws.Range("E49").Value = '  +1.76%  '
$ws.Range("D50").Value = '7.09'
$ws.Range("E50").Value = '  -1.67%  '
$ws.Range("D51").Value = '2.255.32'
$ws.Range("E51").Value = '  +2.73%  '

foreach ($cell in $textCells) {
    $ws.Range($cell).Style = "Normal"
}

